# Applies the "update Sheets via scheduled runner" edits to the Masamune_Profits workbook.
# Generated from the authoritative XML diff: for each touched row, writes the new
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ(HQ) values (columns H-N) cell by cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 45676
$ws.Range("J109").Value = 45676
$ws.Range("L109").Value = 45676
$ws.Range("N109").Value = -48450
$ws.Range("H114").Value = 41212
$ws.Range("J114").Value = 41212
$ws.Range("L114").Value = 41212
$ws.Range("N114").Value = -49890
$ws.Range("H128").Value = 46770
$ws.Range("J128").Value = 46770
$ws.Range("L128").Value = 46770
$ws.Range("N128").Value = -56730
$ws.Range("H130").Value = 37192
$ws.Range("J130").Value = 43990
$ws.Range("L130").Value = 43990
$ws.Range("N130").Value = -54030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 36073
$ws.Range("J107").Value = 36073
$ws.Range("L107").Value = 36073
$ws.Range("N107").Value = -43753
$ws.Range("H121").Value = 42170
$ws.Range("J121").Value = 42170
$ws.Range("L121").Value = 42170
$ws.Range("N121").Value = -45664
$ws.Range("H125").Value = 47851.25
$ws.Range("J125").Value = 47851.25
$ws.Range("L125").Value = 47851.25
$ws.Range("N125").Value = -57691.25
$ws.Range("H128").Value = 47952.668
$ws.Range("J128").Value = 47952.668
$ws.Range("L128").Value = 47952.668
$ws.Range("N128").Value = -57912.668
$ws.Range("H133").Value = 46604
$ws.Range("J133").Value = 46604
$ws.Range("L133").Value = 46604
$ws.Range("N133").Value = -51664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 19875
$ws.Range("J40").Value = 19875
$ws.Range("L40").Value = 19875
$ws.Range("N40").Value = -20405
$ws.Range("H96").Value = 11360.714
$ws.Range("I96").Value = 2721.4285
$ws.Range("K96").Value = 2721.4285
$ws.Range("M96").Value = 24.57150000000001
$ws.Range("H112").Value = 47461
$ws.Range("J112").Value = 47461
$ws.Range("L112").Value = 47461
$ws.Range("N112").Value = -50415
$ws.Range("H119").Value = 47425
$ws.Range("J119").Value = 47425
$ws.Range("L119").Value = 47425
$ws.Range("N119").Value = -57101
$ws.Range("H120").Value = 48753
$ws.Range("J120").Value = 48753
$ws.Range("L120").Value = 48753
$ws.Range("N120").Value = -58429
$ws.Range("H124").Value = 50757.332
$ws.Range("J124").Value = 50757.332
$ws.Range("L124").Value = 50757.332
$ws.Range("N124").Value = -60577.332
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49923.332
$ws.Range("J20").Value = 49923.332
$ws.Range("L20").Value = 49923.332
$ws.Range("N20").Value = -50395.332
$ws.Range("H30").Value = 49923.332
$ws.Range("J30").Value = 49923.332
$ws.Range("L30").Value = 49923.332
$ws.Range("N30").Value = -50105.332
$ws.Range("H99").Value = 1558.7407
$ws.Range("I99").Value = 1326
$ws.Range("K99").Value = 1326
$ws.Range("M99").Value = 172
$ws.Range("H100").Value = 47776
$ws.Range("J100").Value = 47776
$ws.Range("L100").Value = 47776
$ws.Range("N100").Value = -49940
$ws.Range("H110").Value = 38003.332
$ws.Range("J110").Value = 38003.332
$ws.Range("L110").Value = 38003.332
$ws.Range("N110").Value = -46183.332
$ws.Range("H111").Value = 41223
$ws.Range("J111").Value = 41223
$ws.Range("L111").Value = 41223
$ws.Range("N111").Value = -49403
$ws.Range("H116").Value = 48244
$ws.Range("J116").Value = 48244
$ws.Range("L116").Value = 48244
$ws.Range("N116").Value = -57422
$ws.Range("H118").Value = 44734
$ws.Range("J118").Value = 44734
$ws.Range("L118").Value = 44734
$ws.Range("N118").Value = -48048
$ws.Range("H126").Value = 1558.7407
$ws.Range("I126").Value = 1326
$ws.Range("K126").Value = 3978
$ws.Range("M126").Value = -1508
$ws.Range("H128").Value = 49923.332
$ws.Range("J128").Value = 49923.332
$ws.Range("L128").Value = 49923.332
$ws.Range("N128").Value = -59883.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 125354.12
$ws.Range("I125").Value = 377523.75
$ws.Range("J125").Value = 6686.0586
$ws.Range("K125").Value = 1132571.25
$ws.Range("L125").Value = 20058.1758
$ws.Range("M125").Value = -1127651.25
$ws.Range("N125").Value = -29898.1758
$ws.Range("H134").Value = 41713916
$ws.Range("I134").Value = 100104660
$ws.Range("J134").Value = 6244.857
$ws.Range("K134").Value = 300313980
$ws.Range("L134").Value = 18734.571
$ws.Range("M134").Value = -300308910
$ws.Range("N134").Value = -28874.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874
$ws.Range("H114").Value = 48714
$ws.Range("J114").Value = 48714
$ws.Range("L114").Value = 48714
$ws.Range("N114").Value = -57392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H108").Value = 39213
$ws.Range("J108").Value = 39213
$ws.Range("L108").Value = 39213
$ws.Range("N108").Value = -46893
$ws.Range("H111").Value = 43966
$ws.Range("J111").Value = 43966
$ws.Range("L111").Value = 43966
$ws.Range("N111").Value = -52146
$ws.Range("H124").Value = 34641.668
$ws.Range("J124").Value = 34641.668
$ws.Range("L124").Value = 34641.668
$ws.Range("N124").Value = -44461.668
$ws.Range("H127").Value = 50233.2
$ws.Range("J127").Value = 50233.2
$ws.Range("L127").Value = 50233.2
$ws.Range("N127").Value = -60153.2
$ws.Range("H128").Value = 38198.4
$ws.Range("J128").Value = 38198.4
$ws.Range("L128").Value = 38198.4
$ws.Range("N128").Value = -48158.4
$ws.Range("H130").Value = 44406.25
$ws.Range("J130").Value = 43892.855
$ws.Range("L130").Value = 43892.855
$ws.Range("N130").Value = -53932.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 46974.668
$ws.Range("J110").Value = 46974.668
$ws.Range("L110").Value = 46974.668
$ws.Range("N110").Value = -55154.668
$ws.Range("H120").Value = 42196.4
$ws.Range("J120").Value = 42196.4
$ws.Range("L120").Value = 42196.4
$ws.Range("N120").Value = -51872.4
$ws.Range("H123").Value = 36740.145
$ws.Range("J123").Value = 37030.168
$ws.Range("L123").Value = 37030.168
$ws.Range("N123").Value = -46830.168
$ws.Range("H128").Value = 50707
$ws.Range("J128").Value = 50707
$ws.Range("L128").Value = 50707
$ws.Range("N128").Value = -60667
$ws.Range("H131").Value = 49211
$ws.Range("J131").Value = 49211
$ws.Range("L131").Value = 49211
$ws.Range("N131").Value = -59291
$ws.Range("H133").Value = 82868.39999999999
$ws.Range("J133").Value = 82868.39999999999
$ws.Range("L133").Value = 82868.39999999999
$ws.Range("N133").Value = -92988.39999999999
